# Update cryptos list data (rows 2-51) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Bitcoin"
$ws.Range("C2").Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Range("D2").Value = "89.794.39"
$ws.Range("E2").Value = "  -0.91%  "

$ws.Range("B3").Value = "Ethereum"
$ws.Range("C3").Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Range("D3").Value = "3.078.91"
$ws.Range("E3").Value = "  -0.74%  "

$ws.Range("B4").Value = "TetherUSD"
$ws.Range("C4").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.26%  "

$ws.Range("B5").Value = "Solana"
$ws.Range("C5").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D5").Value = "'240.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.70%  "

$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "'616.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.23%  "

$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").Value = "'1.14"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.77%  "

$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").Value = "'0.362"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("B9").Value = "USDC"
$ws.Range("C9").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D9").Value = "'1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.11%  "

$ws.Range("B10").Value = "LidoStakedEther"
$ws.Range("C10").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D10").Value = "3.073.84"
$ws.Range("E10").Value = "  -1.03%  "

$ws.Range("B11").Value = "Cardano"
$ws.Range("C11").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D11").Value = "'0.730"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.24%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.202"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.79%  "

$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").Value = "'0.0000244"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.36%  "

$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "'34.50"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.26%  "

$ws.Range("B15").Value = "Toncoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D15").Value = "'5.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.83%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "89.752.71"
$ws.Range("E16").Value = "  -0.93%  "

$ws.Range("B17").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C17").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D17").Value = "3.645.17"
$ws.Range("E17").Value = "  -1.30%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.071.54"
$ws.Range("E18").Value = "  -1.65%  "

$ws.Range("B19").Value = "SuiNetwork"
$ws.Range("C19").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D19").Value = "'3.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.58%  "

$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'14.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.04%  "

$ws.Range("B21").Value = "PEPE"
$ws.Range("C21").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D21").Value = "'0.0000207"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.45%  "

$ws.Range("B22").Value = "Polkadot"
$ws.Range("C22").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D22").Value = "'5.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.86%  "

$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "'435.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.10%  "

$ws.Range("B24").Value = "Uniswap"
$ws.Range("C24").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D24").Value = "'8.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.72%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'90.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.32%  "

$ws.Range("B26").Value = "NEARProtocol"
$ws.Range("C26").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D26").Value = "'5.57"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.81%  "

$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").Value = "'11.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.24%  "

$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "3.243.14"
$ws.Range("E28").Value = "  -2.18%  "

$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").Value = "'0.244"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +19.33%  "

$ws.Range("B31").Value = "Cronos"
$ws.Range("C31").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D31").Value = "'0.174"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.36%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "'0.119"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +33.16%  "

$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'9.07"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.06%  "

$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").Value = "'1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +12.08%  "

$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "'0.165"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +10.24%  "

$ws.Range("B36").Value = "MantraDAO"
$ws.Range("C36").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D36").Value = "'4.30"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +27.89%  "

$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D37").Value = "'7.57"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.69%  "

$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").Value = "'26.08"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.20%  "

$ws.Range("B39").Value = "PancakeSwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D39").Value = "'1.89"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.09%  "

$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "'482.34"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.16%  "

$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "'3.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.91%  "

$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D42").Value = "'1.27"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.46%  "

$ws.Range("B43").Value = "PolygonEcosystemToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D43").Value = "'0.415"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.12%  "

$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D44").Value = "'22.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.18%  "

$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").Value = "'1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.03%  "

$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").Value = "'154.47"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.31%  "

$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "'1.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.19%  "

$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "'0.679"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.84%  "

$ws.Range("B49").Value = "ImmutableX"
$ws.Range("C49").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D49").Value = "'1.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.12%  "

$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D50").Value = "'44.06"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.32%  "

$ws.Range("B51").Value = "FirstDigitalUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D51").Value = "'0.999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.32%  "
